# "My profile test added"
# Adds a new "emails" worksheet with sample invalid-email test data,
# trims the "users" worksheet login list down to its last 3 entries,
# and updates the active/selected tab state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Create the new "emails" worksheet after the last existing sheet
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$emails = $wb.Worksheets.Add($null, $lastSheet)
$emails.Name = "emails"

$emails.Cells.Item(1, 1).Value = "invalid emails"
$emails.Cells.Item(2, 1).Value = "plainaddress"
$emails.Cells.Item(3, 1).Value = "#@%^%#$@#$@#.com"
$emails.Cells.Item(4, 1).Value = "@domain.com"
$emails.Cells.Item(5, 1).Value = "Joe Smith <email@domain.com>"
$emails.Cells.Item(6, 1).Value = "email.domain.com"
$emails.Cells.Item(7, 1).Value = "email@domain@domain.com"
$emails.Cells.Item(8, 1).Value = ".email@domain.com"
$emails.Cells.Item(9, 1).Value = "email.@domain.com"
$emails.Cells.Item(10, 1).Value = "email..email@domain.com"
$emails.Cells.Item(11, 1).Value = "あいうえお@domain.com"
$emails.Cells.Item(12, 1).Value = "email@domain.com (Joe Smith)"
$emails.Cells.Item(13, 1).Value = "email@domain"
$emails.Cells.Item(14, 1).Value = "email@-domain.com"
$emails.Cells.Item(15, 1).Value = "email@domain.web"
$emails.Cells.Item(16, 1).Value = "email@111.222.333.44444"
$emails.Cells.Item(17, 1).Value = "email@domain..com"

# Style the sample rows (not the header, not the hyperlinked row) with
# Segoe UI / dark grey text, the same look used elsewhere for this data.
$top = $emails.Range("A2:A8")
$top.Font.Name = "Segoe UI"
$top.Font.Color = 3355443

$bottom = $emails.Range("A10:A17")
$bottom.Font.Name = "Segoe UI"
$bottom.Font.Color = 3355443

# Taller rows for the styled entries (matches the "bigger" look of the list)
$emails.Rows.Item(2).RowHeight = 16.5
$emails.Rows.Item(3).RowHeight = 16.5
$emails.Rows.Item(4).RowHeight = 16.5
$emails.Rows.Item(5).RowHeight = 16.5
$emails.Rows.Item(6).RowHeight = 16.5
$emails.Rows.Item(7).RowHeight = 16.5
$emails.Rows.Item(8).RowHeight = 16.5
$emails.Rows.Item(10).RowHeight = 16.5
$emails.Rows.Item(11).RowHeight = 16.5
$emails.Rows.Item(12).RowHeight = 16.5
$emails.Rows.Item(13).RowHeight = 16.5
$emails.Rows.Item(14).RowHeight = 16.5
$emails.Rows.Item(15).RowHeight = 16.5
$emails.Rows.Item(16).RowHeight = 16.5
$emails.Rows.Item(17).RowHeight = 16.5

# The "email.@domain.com" sample is a live mailto hyperlink
$emails.Hyperlinks.Add($emails.Range("A9"), "mailto:email.@domain.com")

$emails.Columns.Item(1).ColumnWidth = 32.42578125

$emails.Range("A17").Select()
$emails.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 2. Trim the "users" worksheet: keep only the last 3 login rows
# ---------------------------------------------------------------
$users = $wb.Worksheets.Item("users")
$users.Range("8:18").EntireRow.Delete()

$users.Range("A5").Value = "jypudiry@amail.club"
$users.Range("B5").Value = "fEpkeCRVzY"
$users.Range("A6").Value = "dybo@banit.me"
$users.Range("B6").Value = "cqqCMyVmsi"
$users.Range("A7").Value = "buguma@duck2.club"
$users.Range("B7").Value = "tx5clQNdOu"

# ---------------------------------------------------------------
# 3. Fix up tab selection: "users" becomes the active tab
# ---------------------------------------------------------------
$users.Activate()
